# Applies the task-sheet update described by the commit:
#  - Adds two new rows (40 & 41) describing the new "Cond Temp Difference"
#    and "Cond Ambient Subcooling Difference" tasks to the bottom of the
#    taskSheet.
#  - Updates the frozen-pane top-left cell / selection to match the new
#    view state recorded for the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("taskSheet")

# Copy the formatting (styles / row height) of the last existing data row
# down onto the two new rows before filling in the values.
$ws.Range("A39:H39").Copy() | Out-Null
$ws.Range("A40:H41").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Rows.Item(40).RowHeight = 45
$ws.Rows.Item(41).RowHeight = 45

# --- Fill in the new shared strings in the same order they were authored ---
$ws.Range("A40").Value = 'Finds the difference between the cond variables'
$ws.Range("C40").Value = '25'
$ws.Range("B41").Value = 'Cond Ambient Subcooling Difference'
$ws.Range("B40").Value = 'Cond Temp Difference'
$ws.Range("E40").Value = 'Cond Air TD `%rackname`'
$ws.Range("D40").Value = 'Cond Outlet Air Temperature `%rackname`,Cond Inlet Air Temperature `%rackname`'
$ws.Range("D41").Value = 'Cond Downleg Temperature `%rackname`,Cond Inlet Air Temperature `%rackname`'
$ws.Range("E41").Value = 'Cond Ambient Subcooling `%rackname`'

# --- Remaining cells reuse already-existing shared strings -----------------
$ws.Range("A41").Value = 'Finds the difference between the cond variables'
$ws.Range("C41").Value = '25'
$ws.Range("F40").Value = '* * * * *'
$ws.Range("G40").Value = '0'
$ws.Range("H40").Value = '1'
$ws.Range("F41").Value = '* * * * *'
$ws.Range("G41").Value = '0'
$ws.Range("H41").Value = '1'

# --- Update the view state (frozen pane / active selection) ------------
$win = $excel.ActiveWindow
$win.ScrollRow = 2
$ws.Range("E22").Select() | Out-Null
